$wb = $excel.ActiveWorkbook
$deployments = $wb.Worksheets.Item("Deployments")
$dbbackups   = $wb.Worksheets.Item("DBBackups")

# --- 1. Add a brand-new, empty "Sheet1" right before DBBackups. ---------
#     (created first so it gets the lower new sheetId = 3)
$sheet1 = $wb.Worksheets.Add($dbbackups)
$sheet1.Name = "Sheet1"
$sheet1.Columns.Item(1).ColumnWidth = 42

# --- 2. Duplicate "Deployments" -> becomes "Deployments_9th" -----------
#     (created second so it gets the next new sheetId = 4); Excel's
#     Worksheet.Copy places the duplicate immediately before the sheet
#     it was copied from, i.e. right at the front of the workbook.
$deployments.Copy($deployments)
$depl9th = $wb.Worksheets.Item(1)
$depl9th.Name = "Deployments_9th"

# --- 3. Update the new sheet's data with the additional deployments ----
$depl9th.Range("A2").Value = "\\vrgefs01\shared\IT\Programme Victory\Releases\Web Apps\20180901\VRUKL.6817.28984"
$depl9th.Range("A3").Value = "\\vrgefs01\shared\IT\Programme Victory\Releases\Web Apps\20180918\VRUKL.6830.25880"
$depl9th.Range("A4").Value = "\\vrgefs01\shared\IT\Programme Victory\Releases\Web Apps\20180924\VRUKL.6838.24144"
$depl9th.Range("A5").Value = "\\vrgefs01\shared\IT\Programme Victory\Releases\Web Apps\20181002\VRUKL.6849.22911\VRUKL.6849.22911"
$depl9th.Range("A6").Value = "\\vrgefs01\shared\IT\Programme Victory\Releases\Web Apps\20180927\VRUKL.6844.22422"
$depl9th.Range("A7").Value = "\\vrgefs01\shared\IT\Programme Victory\Releases\Web Apps\20181001\VRUKL.6848.24666"

# the copied sheet inherited hyperlinks on A2 and A12 - the new sheet has none
$depl9th.Range("A2:A12").Hyperlinks.Delete() | Out-Null

# clear out the rows that used to hold deployments 8-11 and the old footer row
$depl9th.Range("A8:B11").ClearContents()
$depl9th.Range("A12:B12").ClearContents()

# the old sheet had 2 trailing text rows (13 & 14) that are no longer present
$depl9th.Rows("13:14").Delete()

# match the saved selection/active cell
$depl9th.Range("B8").Select() | Out-Null
$depl9th.Activate() | Out-Null
